$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14:45 down to 15:46
$ws.Rows("14").Insert()

# Populate the new row 14 with the new weekly record (same constant columns as
# the surrounding rows, new Fecha/Volumen values)
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C14").Value = 'Metropolitana'
$ws.Range("D14").Value = 44459
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100114002
$ws.Range("G14").Value = 'Camote'
$ws.Range("H14").Value = 'Sin especificar'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 1060
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 12500
$ws.Range("N14").Value = '$/malla 18 kilos'
$ws.Range("O14").Value = 'Perú'
$ws.Range("P14").Value = 694
$ws.Range("Q14").Value = 18
$ws.Range("R14").Value = 'Hortaliza'
